# Reorders the names listed in the "Recorded By" column (column G) of the
# "Session Analysis Results" sheet. Each cell contains a comma-separated
# list of recorder names/emails; this normalizes their order according to
# a fixed priority so that "system" / "dnasr281@gmail.com" sort first and
# "backup@backdoor.com" sorts last, e.g.:
#   "System, dnasr281@gmail.com"              -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System"              -> "System, backup@backdoor.com"
#   "backup@backdoor.com, System, system"      -> "system, System, backup@backdoor.com"
#   "admin@admin.com, dnasr281@gmail.com"      -> "dnasr281@gmail.com, admin@admin.com"
#
# Note: comparisons must be case-sensitive because "System" and "system"
# are distinct values that sort differently ("system" sorts before
# "System"). PowerShell's default hashtable/switch/-eq are case
# INsensitive, so we use the .Equals() string method (case-sensitive by
# default) to rank each entry, and a manual sort instead of Sort-Object.

function Get-RecorderPriority($name) {
    if ($name.Equals("system")) { return 0 }
    if ($name.Equals("dnasr281@gmail.com")) { return 1 }
    if ($name.Equals("System")) { return 2 }
    if ($name.Equals("admin@admin.com")) { return 3 }
    if ($name.Equals("backup@backdoor.com")) { return 4 }
    return 999
}

function Sort-RecordedByValue($value) {
    $parts = @($value -split ",\s*" | ForEach-Object { $_.Trim() })
    $n = $parts.Count
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt ($n - $i - 1); $j++) {
            $p1 = Get-RecorderPriority $parts[$j]
            $p2 = Get-RecorderPriority $parts[$j + 1]
            if ($p1 -gt $p2) {
                $tmp = $parts[$j]
                $parts[$j] = $parts[$j + 1]
                $parts[$j + 1] = $tmp
            }
        }
    }
    return [string]::Join(", ", $parts)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    $value = $cell.Value2
    if ([string]::IsNullOrEmpty($value)) { continue }
    if ($value.IndexOf(",") -lt 0) { continue }

    $newValue = Sort-RecordedByValue $value

    if (-not $newValue.Equals($value)) {
        $cell.Value2 = $newValue
    }
}
